# chore: update Sheets via scheduled runner
# Applies scraped market-board price/profit refresh values to the leve
# profit tracking sheets (ALC, BSM, CUL, GSM, WVR).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 92 - values refreshed, LeveProfitHQ column no longer populated,
# LeveProfitNQ newly populated instead.
$ws.Range("H92").Value = 1000
$ws.Range("I92").Value = 1000
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1000
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("M92").Value = 248

# Rows 125-141 - newly populated pricing columns (H:N) for leves that
# previously had no market data.
$ws.Range("H125").Value = 464.66666
$ws.Range("I125").Value = 497
$ws.Range("J125").Value = 400
$ws.Range("K125").Value = 4473
$ws.Range("L125").Value = 3600
$ws.Range("M125").Value = -2013
$ws.Range("N125").Value = -8520

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 5866.6665
$ws.Range("I129").Value = 600
$ws.Range("J129").Value = 8500
$ws.Range("K129").Value = 1800
$ws.Range("L129").Value = 25500
$ws.Range("M129").Value = 3200
$ws.Range("N129").Value = -35500

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws.Range("H131").Value = 2287.25
$ws.Range("I131").Value = 2287.25
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 6861.75
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1821.75

$ws.Range("H132").Value = 4949.9546
$ws.Range("I132").Value = 4744.95
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 14234.85
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -11704.85
$ws.Range("N132").Value = -26060

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws.Range("H135").Value = 2113
$ws.Range("I135").Value = 2139.625
$ws.Range("J135").Value = 1900
$ws.Range("K135").Value = 19256.625
$ws.Range("L135").Value = 17100
$ws.Range("M135").Value = -16721.625
$ws.Range("N135").Value = -22170

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H138").Value = 3000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

$ws.Range("H139").Value = 100780
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 100780
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 100780
$ws.Range("N139").Value = -111060

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("H141").Value = 33795
$ws.Range("I141").Value = 33795
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 101385
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -96205

# ---------------------------------------------------------------------
# BSM sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H29").Value = 11477.833
$ws.Range("I29").Value = 12492
$ws.Range("J29").Value = 9449.5
$ws.Range("K29").Value = 12492
$ws.Range("L29").Value = 9449.5
$ws.Range("M29").Value = -12203
$ws.Range("N29").Value = -10027.5

$ws.Range("H36").Value = 10497.75
$ws.Range("I36").Value = 9997
$ws.Range("K36").Value = 9997
$ws.Range("M36").Value = -9463

# ---------------------------------------------------------------------
# CUL sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

# ---------------------------------------------------------------------
# GSM sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0

$ws.Range("H127").Value = 15000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 15000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 15000
$ws.Range("N127").Value = -24920

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0

$ws.Range("H130").Value = 120000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 120000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws.Range("H132").Value = 3146.5
$ws.Range("I132").Value = 3172.3333
$ws.Range("J132").Value = 2914
$ws.Range("K132").Value = 9516.999899999999
$ws.Range("L132").Value = 8742
$ws.Range("M132").Value = -6986.999899999999
$ws.Range("N132").Value = -13802

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws.Range("H134").Value = 54409
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 54409
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 163227
$ws.Range("N134").Value = -168297

$ws.Range("H135").Value = 75000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 75000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws.Range("H136").Value = 44448.5
$ws.Range("I136").Value = 65296
$ws.Range("J136").Value = 37499.332
$ws.Range("K136").Value = 195888
$ws.Range("L136").Value = 112497.996
$ws.Range("M136").Value = -193338
$ws.Range("N136").Value = -117597.996

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("H141").Value = 78999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 78999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 78999
$ws.Range("N141").Value = -89359

# ---------------------------------------------------------------------
# WVR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H58").Value = 34148.5
$ws.Range("I58").Value = 28166.666
$ws.Range("K58").Value = 28166.666
$ws.Range("M58").Value = -27858.666

Write-Host "Applied Rafflesia_Profits market data refresh"
